$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture current (pre-edit) row 2 and row 3 values so we can swap them.
$a2 = $ws.Range("A2").Value2
$b2 = $ws.Range("B2").Value2
$c2 = $ws.Range("C2").Value2
$d2 = $ws.Range("D2").Value2
$e2 = $ws.Range("E2").Value2

$a3 = $ws.Range("A3").Value2
$b3 = $ws.Range("B3").Value2
$c3 = $ws.Range("C3").Value2
$d3 = $ws.Range("D3").Value2
$e3 = $ws.Range("E3").Value2

$e4 = $ws.Range("E4").Value2

# Remove all existing hyperlinks; they will be re-created below with the
# swapped addresses/display text so row/cell mapping stays correct.
$ws.Hyperlinks.Delete()

# Write swapped values: row 2 gets the old row 3 content, row 3 gets the old row 2 content.
$ws.Range("A2").Value = $a3
$ws.Range("B2").Value = $b3
$ws.Range("C2").Value = $c3
$ws.Range("D2").Value = $d3

$ws.Range("A3").Value = $a2
$ws.Range("B3").Value = $b2
$ws.Range("C3").Value = $c2
$ws.Range("D3").Value = $d2

# Re-create the hyperlinks for E2, E3 (swapped) and E4 (unchanged).
$ws.Hyperlinks.Add($ws.Range("E2"), $e3, [Type]::Missing, [Type]::Missing, $e3)
$ws.Hyperlinks.Add($ws.Range("E3"), $e2, [Type]::Missing, [Type]::Missing, $e2)
$ws.Hyperlinks.Add($ws.Range("E4"), $e4, [Type]::Missing, [Type]::Missing, $e4)
